# Add ifoCAST full series evaluation: extend header row with Q8/Q9 and
# populate the new diagonal of forecast-error values for rows 19-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new quarter labels, matching the
# existing header style (bold, centered, bordered) by copying format from I1.
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("I1").Copy($ws.Range("K1"))
$ws.Range("J1").Value = "Q8"
$ws.Range("K1").Value = "Q9"

# --- Row 19
$ws.Range("H19").Value = 2.329873611643036
$ws.Range("I19").Value = 9.174120406567614
$ws.Range("J19").Value = -8.376533676533697
$ws.Range("K19").Value = -0.66713707687755

# --- Row 20
$ws.Range("G20").Value = 2.329873611643053
$ws.Range("H20").Value = 9.174120406567631
$ws.Range("I20").Value = -8.376533676533679
$ws.Range("J20").Value = -0.6671370768775329

# --- Row 21
$ws.Range("F21").Value = 2.379873611643061
$ws.Range("G21").Value = 9.224120406567639
$ws.Range("H21").Value = -8.326533676533671
$ws.Range("I21").Value = -0.6171370768775243

# --- Row 22
$ws.Range("E22").Value = 2.329873611643053
$ws.Range("F22").Value = 9.174120406567631
$ws.Range("G22").Value = -8.376533676533679
$ws.Range("H22").Value = -0.6671370768775329

# --- Row 23
$ws.Range("D23").Value = 2.279873611643053
$ws.Range("E23").Value = 9.223120406567634
$ws.Range("F23").Value = -8.327533676533676
$ws.Range("G23").Value = -0.6181370768775289
$ws.Range("H23").Value = 0.9647755840093826
$ws.Range("I23").Value = -2.001511408339873
$ws.Range("J23").Value = 0.2625197463245897
$ws.Range("K23").Value = -0.1886344086021489

# --- Row 24
$ws.Range("C24").Value = 2.269873611643048
$ws.Range("D24").Value = 9.164120406567623
$ws.Range("E24").Value = -8.346533676533667
$ws.Range("F24").Value = -0.6271370768775294
$ws.Range("G24").Value = 1.015775584009385
$ws.Range("H24").Value = -1.950511408339871
$ws.Range("I24").Value = 0.2935197463245814
$ws.Range("J24").Value = -0.1976344086021493
